$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.203.49"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "2.917.82"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "365.14"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("D6").Value = "103.59"
$ws.Range("E6").Value = "  -6.08%  "
$ws.Range("D7").Value = "0.541"
$ws.Range("E7").Value = "  -4.78%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -6.79%  "
$ws.Range("D10").Value = "37.02"
$ws.Range("E10").Value = "  -5.26%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").Value = "0.0834"
$ws.Range("E12").Value = "  -4.24%  "
$ws.Range("D13").Value = "18.50"
$ws.Range("E13").Value = "  -5.74%  "
$ws.Range("D14").Value = "3.375.59"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("E15").Value = "  -5.18%  "
$ws.Range("D16").Value = "2.913.23"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("E17").Value = "  -3.47%  "
$ws.Range("D18").Value = "51.185.30"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").Value = "3.32"
$ws.Range("E19").Value = "  -6.72%  "
$ws.Range("D20").Value = "7.26"
$ws.Range("E20").Value = "  -4.26%  "
$ws.Range("E21").Value = "  -6.47%  "
$ws.Range("D22").Value = "0.0₃0947"
$ws.Range("E22").Value = "  -3.70%  "
$ws.Range("D23").Value = "68.15"
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("D24").Value = "260.06"
$ws.Range("E24").Value = "  -3.26%  "
$ws.Range("E25").Value = "  -4.30%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "4.33"
$ws.Range("E26").Value = "  +3.65%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.175"
$ws.Range("E27").Value = "  -4.98%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "25.97"
$ws.Range("E29").Value = "  -3.87%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "7.31"
$ws.Range("E30").Value = "  -5.10%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.105"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").Value = "6.16"
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "9.96"
$ws.Range("E33").Value = "  -5.26%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "2.14"
$ws.Range("E34").Value = "  -3.33%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "35.26"
$ws.Range("E35").Value = "  -6.07%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "50.69"
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0423"
$ws.Range("E38").Value = "  -4.87%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "2.82"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "3.14"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "16.99"
$ws.Range("E41").Value = "  -7.37%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "1.87"
$ws.Range("E42").Value = "  -6.68%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "0.114"
$ws.Range("E43").Value = "  -5.10%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "22.66"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "117.90"
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "2.11"
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.068.48"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "3.21"
$ws.Range("E48").Value = "  -7.47%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "2.27"
$ws.Range("E49").Value = "  -8.57%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "3.206.30"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "0.235"
$ws.Range("E51").Value = "  -5.90%  "
